# API Usuario SQL V2
# Update the "Nota Final" (column D) values for several rubric rows.
# D4 holds =SUM(D5:D61) and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = 1

$ws.Range("D16").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0

# Reflect the reviewer's cursor position/selection at save time.
$ws.Range("B23").Select()
